$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Education: append ", Certificate of Entrepreneurship and Innovation " to
#    the "Computer Science Minor" line.
# ---------------------------------------------------------------------------
$found = $d.Content.Find.Execute( `
    "Computer Science Minor", $true, $false, $false, $false, $false, $true, 1, $false, `
    "Computer Science Minor, Certificate of Entrepreneurship and Innovation ", 2)

# ---------------------------------------------------------------------------
# 2. Education: GPA 3.93 -> 3.95
# ---------------------------------------------------------------------------
$found = $d.Content.Find.Execute( `
    "GPA: 3.93", $true, $false, $false, $false, $false, $true, 1, $false, `
    "GPA: 3.95", 2)

# ---------------------------------------------------------------------------
# 3. Move the stray "_GoBack" bookmark out of the Programming Languages line
#    (it used to sit between " JavaFX," and " C, Python, MatLab") -- it will
#    be re-added further down, inside the "Designed C code..." bullet.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 4. Experience bullet: "...output digital and analog signals" ->
#    "...output digital and analog control signals", with the "_GoBack"
#    bookmark sitting right before the final word "signals".
# ---------------------------------------------------------------------------
$found = $d.Content.Find.Execute( `
    "Designed C code on TI microcontroller to interpret serial commands and output digital and analog signals", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Designed C code on TI microcontroller to interpret serial commands and output digital and analog control signals", 2)

$r = $d.Content
$found = $r.Find.Execute("analog control signals", $false)
if ($found) {
    $r.Collapse(1)
    $r.MoveStart(1, 15)   # move past "analog control " (15 characters, incl. trailing space)
    $r.Collapse(1)
    $d.Bookmarks.Add("_GoBack", $r)
}

# ---------------------------------------------------------------------------
# 5. Certifications: "git/github" -> "Git/Github"
# ---------------------------------------------------------------------------
$found = $d.Content.Find.Execute( `
    ": LTspice, KiCad, git/github", $true, $false, $false, $false, $false, $true, 1, $false, `
    ": LTspice, KiCad, Git/Github", 2)

# ---------------------------------------------------------------------------
# 6. Awards table: "Fall 2018, Spring 2019" -> "3 Semesters"
# ---------------------------------------------------------------------------
$found = $d.Content.Find.Execute( `
    "Fall 2018, Spring 2019", $true, $false, $false, $false, $false, $true, 1, $false, `
    "3 Semesters", 2)

Write-Output "done"
